# The "names" pool (Sheet1) had its first entry ("hlzcq1o6") consumed.
# Remove row 1 from Sheet1 - this shifts every remaining name up by one
# row (A2->A1, A3->A2, ... ) and shrinks the used range from A1:A483 to
# A1:A482, exactly like the diff shows.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("used")

$ws1.Rows("1:1").Delete()

# Record that name as now "used" by appending a new row to the "used"
# sheet (row 17) with the id, the source filename and the timestamp at
# which it was used.
$ws2.Range("A17").Value = "hlzcq1o6"
$ws2.Range("B17").Value = "ChatGPT Image 2026年1月18日 07_37_14.png"
$ws2.Range("C17").Value = "2026-01-18 07:39:08"
